$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6702.4
$ws.Range("I43").Value = 6733.6
$ws.Range("J43").Value = 6640
$ws.Range("K43").Value = 6733.6
$ws.Range("L43").Value = 6640
$ws.Range("M43").Value = -6664.6
$ws.Range("N43").Value = -6778
$ws.Range("H69").Value = 17999.777
$ws.Range("I69").Value = 8499.5
$ws.Range("J69").Value = 19187.312
$ws.Range("K69").Value = 25498.5
$ws.Range("L69").Value = 57561.936
$ws.Range("M69").Value = -24624.5
$ws.Range("N69").Value = -59309.936
$ws.Range("H72").Value = 17999.777
$ws.Range("I72").Value = 8499.5
$ws.Range("J72").Value = 19187.312
$ws.Range("K72").Value = 76495.5
$ws.Range("L72").Value = 172685.808
$ws.Range("M72").Value = -72127.5
$ws.Range("N72").Value = -181421.808
$ws.Range("H76").Value = 6250.143
$ws.Range("I76").Value = 4662.875
$ws.Range("J76").Value = 8366.5
$ws.Range("K76").Value = 4662.875
$ws.Range("L76").Value = 8366.5
$ws.Range("M76").Value = -4347.875
$ws.Range("N76").Value = -8996.5
$ws.Range("H79").Value = 6250.143
$ws.Range("I79").Value = 4662.875
$ws.Range("J79").Value = 8366.5
$ws.Range("K79").Value = 4662.875
$ws.Range("L79").Value = 8366.5
$ws.Range("M79").Value = -3570.875
$ws.Range("N79").Value = -10550.5
$ws.Range("H80").Value = 1537.5
$ws.Range("I80").Value = 1350
$ws.Range("J80").Value = 1725
$ws.Range("K80").Value = 4050
$ws.Range("L80").Value = 5175
$ws.Range("M80").Value = -3052
$ws.Range("N80").Value = -7171
$ws.Range("H83").Value = 1537.5
$ws.Range("I83").Value = 1350
$ws.Range("J83").Value = 1725
$ws.Range("K83").Value = 12150
$ws.Range("L83").Value = 15525
$ws.Range("M83").Value = -7158
$ws.Range("N83").Value = -25509
$ws.Range("H137").Value = 15199.953
$ws.Range("I137").Value = 10617.226
$ws.Range("J137").Value = 19504.94
$ws.Range("K137").Value = 31851.678
$ws.Range("L137").Value = 58514.81999999999
$ws.Range("M137").Value = -29301.678
$ws.Range("N137").Value = -63614.81999999999
$ws.Range("H138").Value = 13668.029
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 13668.029
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41004.087
$ws.Range("N138").Value = -51284.087

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5248.7646
$ws.Range("I32").Value = 3998.8572
$ws.Range("J32").Value = 20997.6
$ws.Range("K32").Value = 3998.8572
$ws.Range("L32").Value = 20997.6
$ws.Range("M32").Value = -3711.8572
$ws.Range("N32").Value = -21571.6
$ws.Range("H61").Value = 6310.5
$ws.Range("I61").Value = 4051.0667
$ws.Range("J61").Value = 13088.8
$ws.Range("K61").Value = 4051.0667
$ws.Range("L61").Value = 13088.8
$ws.Range("M61").Value = -3839.0667
$ws.Range("N61").Value = -13512.8
$ws.Range("H74").Value = 4082.2163
$ws.Range("I74").Value = 2002.9166
$ws.Range("J74").Value = 7920.923
$ws.Range("K74").Value = 2002.9166
$ws.Range("L74").Value = 7920.923
$ws.Range("M74").Value = -1128.9166
$ws.Range("N74").Value = -9668.922999999999
$ws.Range("H77").Value = 4082.2163
$ws.Range("I77").Value = 2002.9166
$ws.Range("J77").Value = 7920.923
$ws.Range("K77").Value = 10014.583
$ws.Range("L77").Value = 39604.615
$ws.Range("M77").Value = -5646.583000000001
$ws.Range("N77").Value = -48340.615
$ws.Range("H133").Value = 88998
$ws.Range("J133").Value = 88998
$ws.Range("L133").Value = 88998
$ws.Range("N133").Value = -94058
$ws.Range("H136").Value = 6310.5
$ws.Range("I136").Value = 4051.0667
$ws.Range("J136").Value = 13088.8
$ws.Range("K136").Value = 12153.2001
$ws.Range("L136").Value = 39266.39999999999
$ws.Range("M136").Value = -9603.2001
$ws.Range("N136").Value = -44366.39999999999
$ws.Range("H139").Value = 150000
$ws.Range("J139").Value = 150000
$ws.Range("L139").Value = 150000
$ws.Range("N139").Value = -160280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2458.5652
$ws.Range("I134").Value = 1424.3889
$ws.Range("K134").Value = 4273.1667
$ws.Range("M134").Value = -1738.1667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8202.878000000001
$ws.Range("I31").Value = 4704.4
$ws.Range("J31").Value = 9331.419
$ws.Range("K31").Value = 4704.4
$ws.Range("L31").Value = 9331.419
$ws.Range("M31").Value = -4409.4
$ws.Range("N31").Value = -9921.419
$ws.Range("H34").Value = 8202.878000000001
$ws.Range("I34").Value = 4704.4
$ws.Range("J34").Value = 9331.419
$ws.Range("K34").Value = 4704.4
$ws.Range("L34").Value = 9331.419
$ws.Range("M34").Value = -4502.4
$ws.Range("N34").Value = -9735.419

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.96
$ws.Range("I2").Value = 34.3125
$ws.Range("J2").Value = 36.11111
$ws.Range("K2").Value = 205.875
$ws.Range("L2").Value = 216.66666
$ws.Range("M2").Value = -92.875
$ws.Range("N2").Value = -442.66666
$ws.Range("H17").Value = 1225.8889
$ws.Range("I17").Value = 1484.6666
$ws.Range("J17").Value = 708.3333
$ws.Range("K17").Value = 4453.9998
$ws.Range("L17").Value = 2124.9999
$ws.Range("M17").Value = -4284.9998
$ws.Range("N17").Value = -2462.9999
$ws.Range("H18").Value = 266.46667
$ws.Range("I18").Value = 174.61539
$ws.Range("J18").Value = 863.5
$ws.Range("K18").Value = 523.84617
$ws.Range("L18").Value = 2590.5
$ws.Range("M18").Value = -354.84617
$ws.Range("N18").Value = -2928.5
$ws.Range("H34").Value = 4892.8335
$ws.Range("I34").Value = 279
$ws.Range("J34").Value = 7199.75
$ws.Range("K34").Value = 837
$ws.Range("L34").Value = 21599.25
$ws.Range("M34").Value = -753
$ws.Range("N34").Value = -21767.25
$ws.Range("H39").Value = 9437.5
$ws.Range("J39").Value = 9437.5
$ws.Range("L39").Value = 28312.5
$ws.Range("N39").Value = -28900.5
$ws.Range("H46").Value = 1231.9333
$ws.Range("I46").Value = 783.8570999999999
$ws.Range("J46").Value = 1624
$ws.Range("K46").Value = 2351.5713
$ws.Range("L46").Value = 4872
$ws.Range("M46").Value = -2260.5713
$ws.Range("N46").Value = -5054
$ws.Range("H55").Value = 9362.25
$ws.Range("J55").Value = 12333
$ws.Range("L55").Value = 36999
$ws.Range("N55").Value = -37353
$ws.Range("H68").Value = 5031.8965
$ws.Range("J68").Value = 6020.2104
$ws.Range("L68").Value = 18060.6312
$ws.Range("N68").Value = -19682.6312
$ws.Range("H71").Value = 5031.8965
$ws.Range("J71").Value = 6020.2104
$ws.Range("L71").Value = 54181.8936
$ws.Range("N71").Value = -62293.8936
$ws.Range("H121").Value = 1266263.4
$ws.Range("I121").Value = 426.125
$ws.Range("J121").Value = 6329612.5
$ws.Range("K121").Value = 1278.375
$ws.Range("L121").Value = 18988837.5
$ws.Range("M121").Value = 31.625
$ws.Range("N121").Value = -18991457.5
$ws.Range("H131").Value = 758768.4399999999
$ws.Range("I131").Value = 920.1667
$ws.Range("J131").Value = 3032313.2
$ws.Range("K131").Value = 2760.5001
$ws.Range("L131").Value = 9096939.600000001
$ws.Range("M131").Value = 2279.4999
$ws.Range("N131").Value = -9107019.600000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3393.8286
$ws.Range("I122").Value = 2733.6553
$ws.Range("K122").Value = 8200.965899999999
$ws.Range("M122").Value = -5750.965899999999
$ws.Range("H132").Value = 4765.375
$ws.Range("I132").Value = 4223.657
$ws.Range("J132").Value = 8557.4
$ws.Range("K132").Value = 12670.971
$ws.Range("L132").Value = 25672.2
$ws.Range("M132").Value = -10140.971
$ws.Range("N132").Value = -30732.2

# --- Remove cells that no longer exist after the update ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M138").ClearContents()
